# Loca_Keys_Strings_DE_EN.xlsx - add new localized strings for the
# "get free chips" modal, and clean up the old lost-connection row that had
# a stray whitespace-only placeholder cell below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the stray whitespace-only placeholder that used to sit in B72
#     (it was the only reference to that shared string, so removing it
#     also drops the now-unused shared string entry). Row 71 itself
#     (key + EN/DE values) is left untouched. ---
$ws.Range("B72").ClearContents()
$ws.Rows(71).RowHeight = 60

# --- Prepare formatting for the three new rows by copying from existing
#     rows that already carry the right look (fill + wrap combinations),
#     so no redundant style entries get created. ---

# Row 72 pattern: A filled, B filled+wrap, C filled (no wrap)
$ws.Range("A71:B71").Copy()
$ws.Range("A72:B72").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("C72").PasteSpecial(-4122)

# Row 73 pattern: A filled, B filled+wrap, C filled+wrap
$ws.Range("A14:C14").Copy()
$ws.Range("A73:C73").PasteSpecial(-4122)

# Row 74 pattern: A filled, B filled, C filled (no wrap)
$ws.Range("A2:C2").Copy()
$ws.Range("A74:C74").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Fill in the new row values (keys first, then the short EN/DE pairs,
#     then the long content pair last - mirrors how this was authored). ---

$ws.Range("A72").Value = "global_get-free-chips-modal_header"
$ws.Range("A73").Value = "global_get-free-chips-modal_content"
$ws.Range("A74").Value = "global_get-free-chips-modal_btn-txt"

$ws.Range("B72").Value = "Refuel  "
$ws.Range("B74").Value = "Get Your Free Chips"

$ws.Range("C72").Value = "Auftanken"
$ws.Range("C74").Value = "Gratis Chips Holen"

$ws.Range("C73").Value = "Oh nein, es scheint, als würden Ihnen die Chips ausgehen! Aber keine Sorge, hier ist eine neue Charge Chips für Sie, damit Sie weiterspielen können!"
$ws.Range("B73").Value = "Oh noes, it seems like you're running out of chips! But don't worry, here's a fresh batch of chips for you so you can continue playing!"

$ws.Rows(73).RowHeight = 75

# Column A now holds some longer keys ("global_get-free-chips-modal_*"),
# so widen it to fit like the original author's sheet did.
$ws.Columns.Item(1).AutoFit()

$ws.Range("C70").Select()
